# CaseAndFatalityDemographicsData 2021-07-09 update
# Updates the Case and Fatality demographic counts (Age Group, Gender,
# Race/Ethnicity — for both Cases and Fatalities) to the refreshed totals.
# Percentage columns are formulas (B/$B$total) and recompute automatically.

$wb = $excel.ActiveWorkbook

# ---- Sheet 1: Cases by Age Group ----
$ws = $wb.Worksheets.Item("Cases by Age Group")
$ws.Range("B3").Value  = 1414
$ws.Range("B4").Value  = 3909
$ws.Range("B5").Value  = 15925
$ws.Range("B6").Value  = 17483
$ws.Range("B7").Value  = 15331
$ws.Range("B8").Value  = 12939
$ws.Range("B9").Value  = 4683
$ws.Range("B10").Value = 3169
$ws.Range("B11").Value = 1923
$ws.Range("B12").Value = 1274
$ws.Range("B13").Value = 1966

# ---- Sheet 2: Cases by Gender ----
$ws = $wb.Worksheets.Item("Cases by Gender")
$ws.Range("B2").Value = 27489
$ws.Range("B3").Value = 51925

# ---- Sheet 3: Cases by RaceEthnicity ----
$ws = $wb.Worksheets.Item("Cases by RaceEthnicity")
$ws.Range("B2").Value = 990
$ws.Range("B3").Value = 13172
$ws.Range("B4").Value = 28758
$ws.Range("B5").Value = 600
$ws.Range("B6").Value = 28044
$ws.Range("B7").Value = 8749

# ---- Sheet 4: Fatalities by Age Group ----
$ws = $wb.Worksheets.Item("Fatalities by Age Group")
$ws.Range("B4").Value  = 34
$ws.Range("B5").Value  = 272
$ws.Range("B6").Value  = 901
$ws.Range("B7").Value  = 2620
$ws.Range("B8").Value  = 5896
$ws.Range("B9").Value  = 4871
$ws.Range("B10").Value = 6255
$ws.Range("B11").Value = 6879
$ws.Range("B12").Value = 6770
$ws.Range("B13").Value = 16937

# ---- Sheet 5: Fatalities by Gender ----
$ws = $wb.Worksheets.Item("Fatalities by Gender")
$ws.Range("B2").Value = 21575
$ws.Range("B3").Value = 29879

# ---- Sheet 6: Fatalities by Race-Ethnicity ----
$ws = $wb.Worksheets.Item("Fatalities by Race-Ethnicity")
$ws.Range("B2").Value = 1093
$ws.Range("B3").Value = 5256
$ws.Range("B4").Value = 23887
$ws.Range("B5").Value = 284
$ws.Range("B6").Value = 20912
